$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A. This shifts the existing
# Code / Description / Definition columns from A/B/C to B/C/D.
$ws.Range("A1").EntireColumn.Insert()

# New header for the inserted column
$ws.Range("A1").Value = "Version"

# Make sure the "1.0" version values are stored as text (not coerced to
# the number 1) before writing them, then strip the temporary formatting
# back off so the cells end up with no explicit style, matching plain
# text cells elsewhere in the sheet.
$ws.Range("A2:A10").NumberFormat = "@"
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value = "1.0"
}
$ws.Range("A2:A10").ClearFormats()
